{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph directly above it) that used to follow the\n// \"Requisitos\" section, while leaving the final blank / page-break\n// paragraphs at the very end of the document intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4038: ...\" requirement paragraph; the three paragraphs\n// that immediately follow it (blank, \"Ver no Jupiter...\", \"\u00a9 2020 ...\")\n// are the ones being removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4038\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4038' requirements paragraph.\");\n}\n\nconst toDelete = [];\nfor (let offset = 1; offset <= 3; offset++) {\n  const idx = anchorIndex + offset;\n  if (idx >= items.length) break;\n  const text = items[idx].text;\n  if (\n    text.trim() === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"Powered by Jekyll\") !== -1\n  ) {\n    toDelete.push(items[idx]);\n  } else {\n    break;\n  }\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph directly above it) that used to follow the\n# \"Requisitos\" section, while leaving the final blank / page-break\n# paragraphs at the very end of the document intact.\n\n$d = $word.ActiveDocument\n\n# Locate the \"LOQ4038: ...\" requirement paragraph via Find, then map the\n# found range back to its paragraph index.\n$rng = $d.Content\n$rng.Find.Execute(\"LOQ4038\") | Out-Null\n$target = $rng.Start\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $target -and $p.Range.End -ge $target) {\n        $anchorIndex = $i\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOQ4038' requirements paragraph.\"\n}\n\n# The (up to) three paragraphs right after it -- a blank paragraph,\n# \"Ver no Jupiter Salvar em pdf Salvar em docx\" and the \"(c) 2020 ...\"\n# copyright footer -- are the ones being removed.\n$toDelete = @()\nfor ($offset = 1; $offset -le 3; $offset++) {\n    $idx = $anchorIndex + $offset\n    if ($idx -gt $d.Paragraphs.Count) { break }\n    $t = $d.Paragraphs.Item($idx).Range.Text.Trim()\n    if ($t -eq \"\" -or $t -like \"*Ver no Jupiter*\" -or $t -like \"*Powered by Jekyll*\") {\n        $toDelete += $idx\n    } else {\n        break\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\nfor ($i = $toDelete.Count - 1; $i -ge 0; $i--) {\n    $d.Paragraphs.Item($toDelete[$i]).Range.Delete()\n}\n"}
